$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ajout dans tableau de bord : nouvelle ligne pour la contribution d'Anthony
# (meme date que la ligne precedente de Francois).
$ws.Range("B11").Value = 44183

# Reprend le format de date (style) de la cellule B10 au lieu de recreer un
# format numerique personnalise.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C11").Value = "Anthony"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "Partie 1"

$ws.Range("K20").Select() | Out-Null
